$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from its old position (inside the
#    "Plotall sometimes does not work ... skimmed datasets" sentence)
#    to a new position inside "There is an error in the plotting of
#    the prediction ..." (splitting "plotting" into "plott" | "ing").
# -------------------------------------------------------------------

# Remove the existing bookmark first.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-join the sentence that used to be split around the old bookmark
# location into a single run again.
$r = $d.Content
$r.Find.Execute(
    " sometimes does not work on skimmed datasets, but this should be resolved when using exclude",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " sometimes does not work on skimmed datasets, but this should be resolved when using exclude",
    2) | Out-Null

# Find the new split point ("...the plott" | "ing of the prediction...")
# and drop the bookmark there.
$r2 = $d.Content
$r2.Find.Execute("There is an error in the plott", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $r2.End
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos)) | Out-Null

# -------------------------------------------------------------------
# 2) Strike through the "Autowin can't deal with lmer, add centre
#    argument to autowin" bullet (issue resolved).
# -------------------------------------------------------------------

$r3 = $d.Content
$r3.Find.Execute("Autowin can", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $r3.Paragraphs(1)
$para.Range.Font.StrikeThrough = 1
